# Create two new sheets from the "Names" template, mirroring the
# "right click tab -> Move or Copy -> Create a copy" workflow performed
# twice: the first copy is later cleared out (blank placeholder), the
# second copy keeps the template's data and becomes the active sheet.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Names")

# --- First copy: becomes "Names (13)", then gets emptied out ---
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $afterSheet)
$copy1 = $wb.Worksheets.Item($wb.Worksheets.Count)

# Wipe the content/formatting of the duplicated sheet so it ends up
# blank, the same way a user would select everything and delete it.
$copy1.UsedRange.Clear()
$copy1.Cells.UnMerge()
foreach ($hl in $copy1.Hyperlinks) { $hl.Delete() }
$copy1.Rows("1:12").Delete()

# --- Second copy: becomes "Names (14)" and keeps the full data ---
$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $afterSheet2)
$copy2 = $wb.Worksheets.Item($wb.Worksheets.Count)

# The newly created, data-filled sheet is the one the user ends up
# looking at.
$copy2.Activate()
